$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(9894.6200000000008, 9849.31,   283.47000000000003, 284.77,             $false, 0.46, 42613.767083333332, $true),
    @(9896.6,             9894.6200000000008, 282.39,    282.45999999999998, $false, 0.02, 42614.674120370371, $true),
    @(9930.25,            9896.6,    280.62,              281.57,            $false, 0.34, 42615.75273148148,  $true)
)

$startRow = 7
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 7).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($r, 8).Value = $data[7]
}
